$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("F2").Value = 1.45
$ws.Range("G2").Value = 1.82
$ws.Range("H2").Value = 5.1
$ws.Range("I2").Value = 870
$ws.Range("J2").Value = 3.5
$ws.Range("K2").Value = 8.6
$ws.Range("L2").Value = 1.32
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 2.94
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 1.73
$ws.Range("Q2").Value = 1.79
$ws.Range("S2").Value = 2.68
$ws.Range("T2").Value = 1.04
$ws.Range("U2").Value = 1.04
$ws.Range("V2").Value = 1.02
$ws.Range("W2").Value = 2.22

# Row 3
$ws.Range("F3").Value = 1.6
$ws.Range("G3").Value = 2.04
$ws.Range("I3").Value = 7.2
$ws.Range("J3").Value = 3.05
$ws.Range("K3").Value = 5.8
$ws.Range("L3").Value = 1.37
$ws.Range("S3").Value = 3.9
$ws.Range("V3").Value = 1.17
$ws.Range("W3").Value = 1.97

# Row 4
$ws.Range("F4").Value = 1.7
$ws.Range("G4").Value = 1.97
$ws.Range("I4").Value = 870
$ws.Range("J4").Value = 3.05
$ws.Range("K4").Value = 950
$ws.Range("N4").Value = 2.42
$ws.Range("P4").Value = 1.49
$ws.Range("Q4").Value = 2.28
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04
$ws.Range("W4").Value = 2.04

# Row 5
$ws.Range("F5").Value = 7.2
$ws.Range("G5").Value = 7.8
$ws.Range("H5").Value = 1.52
$ws.Range("I5").Value = 1.54
$ws.Range("O5").Value = 1.26
$ws.Range("P5").Value = 2.26
$ws.Range("Q5").Value = 1.78
$ws.Range("R5").Value = 1.46
$ws.Range("U5").Value = 1.99
$ws.Range("V5").Value = 2.86
$ws.Range("Z5").Value = 9
$ws.Range("AC5").Value = 10.5
$ws.Range("AF5").Value = 65
$ws.Range("AG5").Value = 27
$ws.Range("AI5").Value = 32
$ws.Range("AJ5").Value = 220
$ws.Range("AO5").Value = 7.4

# Row 6
$ws.Range("R6").Value = 1.92
$ws.Range("S6").Value = 2.04
$ws.Range("Y6").Value = 65
$ws.Range("Z6").Value = 160
$ws.Range("AA6").Value = 570
$ws.Range("AB6").Value = 14
$ws.Range("AE6").Value = 160

# Row 7
$ws.Range("G7").Value = 3.55
$ws.Range("Q7").Value = 1.73
$ws.Range("U7").Value = 2.46
$ws.Range("W7").Value = 1.39

# Row 8
$ws.Range("F8").Value = 1.46
$ws.Range("G8").Value = 1.48
$ws.Range("H8").Value = 7.8
$ws.Range("I8").Value = 8.4
$ws.Range("L8").Value = 1.36
$ws.Range("N8").Value = 4.3
$ws.Range("Q8").Value = 1.82
$ws.Range("S8").Value = 3.1
$ws.Range("V8").Value = 1.13
$ws.Range("Z8").Value = 75
$ws.Range("AA8").Value = 290
$ws.Range("AK8").Value = 15
$ws.Range("AL8").Value = 36

# Row 9
$ws.Range("G9").Value = 3.3
$ws.Range("I9").Value = 2.24
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 5.8
$ws.Range("W9").Value = 1.43
$ws.Range("Y9").Value = 15
$ws.Range("AB9").Value = 18.5
$ws.Range("AN9").Value = 22

# Row 10
$ws.Range("G10").Value = 2.32
$ws.Range("H10").Value = 3.15
$ws.Range("N10").Value = 6
$ws.Range("P10").Value = 2.68
$ws.Range("R10").Value = 1.68
$ws.Range("S10").Value = 2.38
$ws.Range("T10").Value = 1.52
$ws.Range("W10").Value = 1.76
$ws.Range("AA10").Value = 50
$ws.Range("AJ10").Value = 30

# Row 11
$ws.Range("F11").Value = 2.26
$ws.Range("J11").Value = 3.8
$ws.Range("Q11").Value = 1.72
$ws.Range("S11").Value = 2.76
$ws.Range("U11").Value = 2.5
$ws.Range("AC11").Value = 8.6
$ws.Range("AD11").Value = 14
$ws.Range("AM11").Value = 65

# Row 12
$ws.Range("H12").Value = 17
$ws.Range("I12").Value = 19
$ws.Range("S12").Value = 1.67
$ws.Range("T12").Value = 1.81
$ws.Range("Y12").Value = 1000
$ws.Range("Z12").Value = 240
$ws.Range("AD12").Value = 65
$ws.Range("AE12").Value = 250
$ws.Range("AH12").Value = 36
$ws.Range("AL12").Value = 32
$ws.Range("AN12").Value = 2.5

# Row 13
$ws.Range("O13").Value = 1.23
$ws.Range("S13").Value = 2.74
$ws.Range("T13").Value = 1.61
$ws.Range("AJ13").Value = 50
$ws.Range("AO13").Value = 14.5
